$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'werden',
    'spinnen',
    'münzen',
    'heilen',
    'boxen',
    'arten',
    'bitten',
    'räumen',
    'achten',
    'freuen',
    'suchen',
    'flüchten',
    'sperren',
    'schwören',
    'spielen',
    'sorgen',
    'stehlen',
    'biegen',
    'kehren',
    'ärgern',
    'machen',
    'altern',
    'formen',
    'ehren',
    'quälen',
    'pfeifen',
    'zielen',
    'schenken',
    'bergen',
    'kümmern',
    'decken',
    'liegen',
    'wundern',
    'tropfen',
    'heulen',
    'greifen',
    'sprengen',
    'siegen',
    'wehtun',
    'pflanzen',
    'dienen',
    'klettern',
    'backen',
    'schulden',
    'planen',
    'fischen',
    'hören',
    'fallen',
    'fällen',
    'trauen',
    'bauen',
    'gelten',
    'fahren',
    'schwingen',
    'tollen',
    'lockern',
    'schmecken',
    'zögern',
    'mögen',
    'sichern',
    'reizen',
    'rasen',
    'graben',
    'warnen',
    'drehen',
    'irren',
    'töten',
    'hauen',
    'platzen',
    'saufen',
    'dringen',
    'treiben',
    'rufen',
    'schrecken',
    'fließen',
    'feiern',
    'spüren',
    'heben',
    'bellen',
    'jubeln',
    'stecken',
    'gründen',
    'filmen',
    'geben',
    'scheinen',
    'stammen',
    'seufzen',
    'flehen',
    'wirken',
    'schwächen',
    'knarren',
    'weichen',
    'runden',
    'enden',
    'fangen',
    'zünden',
    'trennen',
    'äußern',
    'erben',
    'ändern',
    'folgen',
    'zeigen',
    'malen',
    'schreiten',
    'bluten',
    'kichern',
    'liefern',
    'werfen',
    'grüßen',
    'streichen',
    'scheitern',
    'kosten',
    'wachsen',
    'loben',
    'sterben',
    'wenden',
    'schlucken',
    'mauern',
    'helfen',
    'führen',
    'betteln',
    'wüten',
    'sinken',
    'brauchen',
    'jagen',
    'lügen',
    'klingen',
    'lesen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}
